$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirements")

# Two more rows of test-case / requirement IDs were added below the existing
# non-functional-requirements block (rows 38-41): R2.5/TC_39 and R2.6/TC_40.
# Column C (TC_ID) is filled in before column A (R_ID) so the shared-string
# table picks up the new strings in the same order the author typed them.

$ws.Range("C42").Value = "TC_39"
$ws.Range("C43").Value = "TC_40"
$ws.Range("A42").Value = "R2.5"
$ws.Range("A43").Value = "R2.6"

# Match the look of the rows above (thin border all around, left/top aligned
# text for column A, left aligned for column C) without touching column B,
# which stays completely empty/untouched for these two rows.
$colA = $ws.Range("A42:A43")
$colA.Borders.LineStyle = 1
$colA.HorizontalAlignment = -4131
$colA.VerticalAlignment = -4160

$colC = $ws.Range("C42:C43")
$colC.Borders.LineStyle = 1

# Update the view state to reflect scrolling down to the newly entered rows.
$ws.Range("A41:A43").Select()
